$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 53, pushing the existing rows 53:81 down to 54:82.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new weekly price record.
$ws.Cells.Item(53, 1).Value = 7
$ws.Cells.Item(53, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(53, 3).Value = "Ñuble"
$ws.Cells.Item(53, 4).Value = 45141
$ws.Cells.Item(53, 5).Value = 16
$ws.Cells.Item(53, 6).Value = 100112026
$ws.Cells.Item(53, 7).Value = "Haba"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 50
$ws.Cells.Item(53, 11).Value = 15000
$ws.Cells.Item(53, 12).Value = 15000
$ws.Cells.Item(53, 13).Value = 15000
$ws.Cells.Item(53, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(53, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(53, 16).Value = 600
$ws.Cells.Item(53, 17).Value = 25
$ws.Cells.Item(53, 18).Value = "Hortaliza"

# Keep the date column formatted the same as the other date cells in column D.
$ws.Cells.Item(53, 4).NumberFormat = $ws.Cells.Item(54, 4).NumberFormat
